$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'34.416.42"
$ws.Range("E2").Value = "  -0.09%  "

# Row 3
$ws.Range("D3").Value = "'1.787.51"
$ws.Range("E3").Value = "  -2.20%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'225.01"
$ws.Range("E5").Value = "  -2.36%  "

# Row 6
$ws.Range("D6").Value = "'0.553"
$ws.Range("E6").Value = "  -3.46%  "

# Row 7
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("D8").Value = "'33.02"
$ws.Range("E8").Value = "  +4.70%  "

# Row 9
$ws.Range("D9").Value = "'0.282"
$ws.Range("E9").Value = "  -2.12%  "

# Row 10
$ws.Range("D10").Value = "'0.0663"
$ws.Range("E10").Value = "  -2.55%  "

# Row 11
$ws.Range("D11").Value = "'0.0932"
$ws.Range("E11").Value = "  -0.15%  "

# Row 12
$ws.Range("D12").Value = "'2.045.82"
$ws.Range("E12").Value = "  -2.13%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.799.43"
$ws.Range("E13").Value = "  -1.65%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.01"
$ws.Range("E14").Value = "  +6.80%  "

# Row 15
$ws.Range("D15").Value = "'0.635"
$ws.Range("E15").Value = "  -3.12%  "

# Row 16
$ws.Range("D16").Value = "'34.454.75"
$ws.Range("E16").Value = "  +0.22%  "

# Row 17
$ws.Range("D17").Value = "'4.25"
$ws.Range("E17").Value = "  -1.57%  "

# Row 18
$ws.Range("D18").Value = "'69.20"
$ws.Range("E18").Value = "  -1.87%  "

# Row 19
$ws.Range("D19").Value = "'255.87"
$ws.Range("E19").Value = "  -0.93%  "

# Row 20
$ws.Range("D20").Value = "'0.0₃0748"
$ws.Range("E20").Value = "  -1.56%  "

# Row 21
$ws.Range("E21").Value = "  +0.06%  "

# Row 22
$ws.Range("D22").Value = "'10.43"
$ws.Range("E22").Value = "  -2.29%  "

# Row 23
$ws.Range("D23").Value = "'4.21"
$ws.Range("E23").Value = "  -3.27%  "

# Row 24
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  -4.30%  "

# Row 25
$ws.Range("D25").Value = "'157.74"
$ws.Range("E25").Value = "  -1.13%  "

# Row 26
$ws.Range("D26").Value = "'16.47"
$ws.Range("E26").Value = "  -1.90%  "

# Row 27
$ws.Range("D27").Value = "'7.04"
$ws.Range("E27").Value = "  -2.05%  "

# Row 28
$ws.Range("E28").Value = "  -3.59%  "

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("E30").Value = "  -2.55%  "

# Row 31
$ws.Range("D31").Value = "'0.0514"
$ws.Range("E31").Value = "  -2.11%  "

# Row 32
$ws.Range("E32").Value = "  -2.12%  "

# Row 33
$ws.Range("E33").Value = "  +0.68%  "

# Row 34
$ws.Range("D34").Value = "'1.90"
$ws.Range("E34").Value = "  +4.53%  "

# Row 35
$ws.Range("D35").Value = "'1.456.19"
$ws.Range("E35").Value = "  -5.43%  "

# Row 36
$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  -2.06%  "

# Row 37
$ws.Range("D37").Value = "'0.631"
$ws.Range("E37").Value = "  -1.31%  "

# Row 39
$ws.Range("D39").Value = "'84.00"
$ws.Range("E39").Value = "  -0.76%  "

# Row 40
$ws.Range("E40").Value = "  +1.45%  "

# Row 41
$ws.Range("D41").Value = "'2.35"
$ws.Range("E41").Value = "  -0.39%  "

# Row 42
$ws.Range("D42").Value = "'0.893"
$ws.Range("E42").Value = "  -2.55%  "

# Row 43
$ws.Range("D43").Value = "'2.08"
$ws.Range("E43").Value = "  -2.66%  "

# Row 44
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.0506"
$ws.Range("E44").Value = "  -4.05%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'12.63"
$ws.Range("E45").Value = "  +3.67%  "

# Row 46
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.06"
$ws.Range("E46").Value = "  -2.37%  "

# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'5.90"
$ws.Range("E47").Value = "  +1.00%  "

# Row 48
$ws.Range("D48").Value = "'1.945.78"
$ws.Range("E48").Value = "  -1.82%  "

# Row 49
$ws.Range("D49").Value = "'0.999"
$ws.Range("E49").Value = "  -0.13%  "

# Row 50
$ws.Range("D50").Value = "'99.26"
$ws.Range("E50").Value = "  -0.07%  "

# Row 51
$ws.Range("D51").Value = "'50.35"
$ws.Range("E51").Value = "  -2.69%  "
